$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(175, 4).Value = 44748
$ws.Cells.Item(176, 4).Value = 44699
$ws.Cells.Item(176, 10).Value = 600
$ws.Cells.Item(176, 11).Value = 7000
$ws.Cells.Item(176, 12).Value = 7000
$ws.Cells.Item(176, 13).Value = 7000
$ws.Cells.Item(176, 16).Value = 1167
$ws.Cells.Item(177, 4).Value = 44385
$ws.Cells.Item(177, 10).Value = 1000
$ws.Cells.Item(178, 4).Value = 44291
$ws.Cells.Item(178, 10).Value = 500
$ws.Cells.Item(178, 11).Value = 8000
$ws.Cells.Item(178, 12).Value = 8000
$ws.Cells.Item(178, 13).Value = 8000
$ws.Cells.Item(178, 16).Value = 1333
$ws.Cells.Item(179, 4).Value = 44711
$ws.Cells.Item(179, 10).Value = 700
$ws.Cells.Item(179, 11).Value = 6000
$ws.Cells.Item(179, 12).Value = 6000
$ws.Cells.Item(179, 13).Value = 6000
$ws.Cells.Item(179, 16).Value = 1000
$ws.Cells.Item(180, 4).Value = 44239
$ws.Cells.Item(180, 10).Value = 300
$ws.Cells.Item(180, 11).Value = 8000
$ws.Cells.Item(180, 12).Value = 8000
$ws.Cells.Item(180, 13).Value = 8000
$ws.Cells.Item(180, 16).Value = 1333
$ws.Cells.Item(181, 4).Value = 44658
$ws.Cells.Item(181, 10).Value = 400
$ws.Cells.Item(181, 11).Value = 7500
$ws.Cells.Item(181, 12).Value = 7500
$ws.Cells.Item(181, 13).Value = 7500
$ws.Cells.Item(181, 16).Value = 1250
$ws.Cells.Item(182, 4).Value = 44348
$ws.Cells.Item(182, 10).Value = 500
$ws.Cells.Item(182, 11).Value = 7000
$ws.Cells.Item(182, 12).Value = 7000
$ws.Cells.Item(182, 13).Value = 7000
$ws.Cells.Item(182, 16).Value = 1167
$ws.Cells.Item(183, 4).Value = 44505
$ws.Cells.Item(183, 10).Value = 400
$ws.Cells.Item(183, 11).Value = 8000
$ws.Cells.Item(183, 12).Value = 8000
$ws.Cells.Item(183, 13).Value = 8000
$ws.Cells.Item(183, 16).Value = 1333
$ws.Cells.Item(184, 4).Value = 44746
$ws.Cells.Item(184, 11).Value = 6500
$ws.Cells.Item(184, 12).Value = 6500
$ws.Cells.Item(184, 13).Value = 6500
$ws.Cells.Item(184, 16).Value = 1083
$ws.Cells.Item(185, 4).Value = 44386
$ws.Cells.Item(185, 10).Value = 600
$ws.Cells.Item(185, 11).Value = 8000
$ws.Cells.Item(185, 12).Value = 8000
$ws.Cells.Item(185, 13).Value = 8000
$ws.Cells.Item(185, 16).Value = 1333
$ws.Cells.Item(186, 4).Value = 44690
$ws.Cells.Item(186, 10).Value = 500
$ws.Cells.Item(186, 11).Value = 7000
$ws.Cells.Item(186, 12).Value = 7000
$ws.Cells.Item(186, 13).Value = 7000
$ws.Cells.Item(186, 16).Value = 1167
$ws.Cells.Item(187, 4).Value = 44307
$ws.Cells.Item(187, 10).Value = 400
$ws.Cells.Item(187, 11).Value = 8000
$ws.Cells.Item(187, 12).Value = 8000
$ws.Cells.Item(187, 13).Value = 8000
$ws.Cells.Item(187, 16).Value = 1333
$ws.Cells.Item(188, 4).Value = 44672
$ws.Cells.Item(188, 11).Value = 9000
$ws.Cells.Item(188, 12).Value = 9000
$ws.Cells.Item(188, 13).Value = 9000
$ws.Cells.Item(188, 14).Value = '$/docena de matas'
$ws.Cells.Item(188, 16).Value = 1500
$ws.Cells.Item(188, 17).Value = 6
$ws.Cells.Item(189, 4).Value = 44344
$ws.Cells.Item(189, 10).Value = 500
$ws.Cells.Item(189, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(189, 16).Value = 583
$ws.Cells.Item(189, 17).Value = 12
$ws.Cells.Item(190, 4).Value = 44433
$ws.Cells.Item(190, 10).Value = 600
$ws.Cells.Item(190, 11).Value = 7000
$ws.Cells.Item(190, 12).Value = 7000
$ws.Cells.Item(190, 13).Value = 7000
$ws.Cells.Item(190, 16).Value = 1167
$ws.Cells.Item(191, 4).Value = 44707
$ws.Cells.Item(191, 10).Value = 700
$ws.Cells.Item(191, 11).Value = 6000
$ws.Cells.Item(191, 12).Value = 6000
$ws.Cells.Item(191, 13).Value = 6000
$ws.Cells.Item(191, 16).Value = 1000
$ws.Cells.Item(192, 4).Value = 44747
$ws.Cells.Item(192, 10).Value = 600
$ws.Cells.Item(192, 12).Value = 7000
$ws.Cells.Item(192, 13).Value = 7000
$ws.Cells.Item(192, 16).Value = 1167
$ws.Cells.Item(193, 4).Value = 44421
$ws.Cells.Item(193, 10).Value = 500
$ws.Cells.Item(193, 11).Value = 8000
$ws.Cells.Item(193, 12).Value = 8000
$ws.Cells.Item(193, 13).Value = 8000
$ws.Cells.Item(193, 16).Value = 1333
$ws.Cells.Item(194, 1).Value = 5
$ws.Cells.Item(194, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(194, 3).Value = 'Maule'
$ws.Cells.Item(194, 4).Value = 44442
$ws.Cells.Item(194, 5).Value = 7
$ws.Cells.Item(194, 6).Value = 100112017
$ws.Cells.Item(194, 7).Value = 'Apio'
$ws.Cells.Item(194, 8).Value = 'Americana (o)'
$ws.Cells.Item(194, 9).Value = 'Primera'
$ws.Cells.Item(194, 10).Value = 700
$ws.Cells.Item(194, 11).Value = 7000
$ws.Cells.Item(194, 12).Value = 8000
$ws.Cells.Item(194, 13).Value = 7571
$ws.Cells.Item(194, 14).Value = '$/docena de matas'
$ws.Cells.Item(194, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(194, 16).Value = 1262
$ws.Cells.Item(194, 17).Value = 6
$ws.Cells.Item(194, 18).Value = 'Hortaliza'
$ws.Cells.Item(195, 1).Value = 5
$ws.Cells.Item(195, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(195, 3).Value = 'Maule'
$ws.Cells.Item(195, 4).Value = 44483
$ws.Cells.Item(195, 5).Value = 7
$ws.Cells.Item(195, 6).Value = 100112017
$ws.Cells.Item(195, 7).Value = 'Apio'
$ws.Cells.Item(195, 8).Value = 'Americana (o)'
$ws.Cells.Item(195, 9).Value = 'Primera'
$ws.Cells.Item(195, 10).Value = 600
$ws.Cells.Item(195, 11).Value = 7000
$ws.Cells.Item(195, 12).Value = 7000
$ws.Cells.Item(195, 13).Value = 7000
$ws.Cells.Item(195, 14).Value = '$/docena de matas'
$ws.Cells.Item(195, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(195, 16).Value = 1167
$ws.Cells.Item(195, 17).Value = 6
$ws.Cells.Item(195, 18).Value = 'Hortaliza'
$ws.Cells.Item(196, 1).Value = 5
$ws.Cells.Item(196, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(196, 3).Value = 'Maule'
$ws.Cells.Item(196, 4).Value = 44663
$ws.Cells.Item(196, 5).Value = 7
$ws.Cells.Item(196, 6).Value = 100112017
$ws.Cells.Item(196, 7).Value = 'Apio'
$ws.Cells.Item(196, 8).Value = 'Americana (o)'
$ws.Cells.Item(196, 9).Value = 'Primera'
$ws.Cells.Item(196, 10).Value = 500
$ws.Cells.Item(196, 11).Value = 9000
$ws.Cells.Item(196, 12).Value = 9000
$ws.Cells.Item(196, 13).Value = 9000
$ws.Cells.Item(196, 14).Value = '$/docena de matas'
$ws.Cells.Item(196, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(196, 16).Value = 1500
$ws.Cells.Item(196, 17).Value = 6
$ws.Cells.Item(196, 18).Value = 'Hortaliza'

$ws.Cells.Item(195, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
